$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Remove the two existing hyperlinks (on B27 / B28) before the row insert
# shifts everything down, then re-add them afterwards at the new location
# (the engine doesn't auto-shift hyperlink anchors on row insert).
$ws.Range("B27").Hyperlinks.Delete() | Out-Null
$ws.Range("B28").Hyperlinks.Delete() | Out-Null

# Insert a new blank row above the "In which country do you currently
# live?" question (currently row 20), shifting every row below it down
# by one.
$ws.Rows.Item(20).Insert()

# Populate the newly-inserted row with the new screener question.
$ws.Range("A20").Value = "Do you currently live in the United States?"
$ws.Range("B20").Value = "Yes"

# Style A20: bold, black text (matches the new font used for this
# question header).
$ws.Range("A20").Font.Bold = $true
$ws.Range("A20").Font.Color = 0

# Style B20: 11pt black text, right aligned (matches the new font used
# for this answer cell).
$ws.Range("B20").Font.Size = 11
$ws.Range("B20").Font.Color = 0
$ws.Range("B20").HorizontalAlignment = -4152

# Re-create the hyperlinks at their shifted locations (B28 / B29).
$ws.Hyperlinks.Add($ws.Range("B28"), "mailto:consent_participant@yopmail.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B29"), "mailto:consent_participant@yopmail.com") | Out-Null

# Match the author's saved selection.
$ws.Range("A19").Select() | Out-Null
